# Update the LR-pairs sheet with refreshed TPM-based expression values for
# the "ECs" sending/target cluster. Raw ligand (G/H) and receptor (M/N)
# expression values for ECs change, and the dependent specificity /
# edge-weight columns (I,J,O,P,Q,R,S,T) are updated to the recomputed
# values that follow from the new TPM numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.749936
$ws.Range("H2").Value = 95.249808
$ws.Range("I2").Value = 0.5302851438878331
$ws.Range("J2").Value = 0.5302851438878331
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 353.9116471018561
$ws.Range("R2").Value = 3185.204823916704
$ws.Range("S2").Value = 0.1375987118501673
$ws.Range("T2").Value = 0.1375987118501673
$ws.Range("G3").Value = 31.749936
$ws.Range("H3").Value = 95.249808
$ws.Range("I3").Value = 0.5302851438878331
$ws.Range("J3").Value = 0.5302851438878331
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 879.012757878192
$ws.Range("R3").Value = 7911.114820903727
$ws.Range("S3").Value = 0.3417548537166183
$ws.Range("T3").Value = 0.3417548537166183
$ws.Range("G4").Value = 31.749936
$ws.Range("H4").Value = 95.249808
$ws.Range("I4").Value = 0.5302851438878331
$ws.Range("J4").Value = 0.5302851438878331
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 130.998891855376
$ws.Range("R4").Value = 1178.990026698384
$ws.Range("S4").Value = 0.05093157832104753
$ws.Range("T4").Value = 0.05093157832104753
$ws.Range("I5").Value = 0.3451699599880819
$ws.Range("J5").Value = 0.3451699599880819
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 230.366003041004
$ws.Range("R5").Value = 2073.294027369036
$ws.Range("S5").Value = 0.08956491127682824
$ws.Range("T5").Value = 0.08956491127682824
$ws.Range("I6").Value = 0.3451699599880819
$ws.Range("J6").Value = 0.3451699599880819
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.2224529775023261
$ws.Range("T6").Value = 0.2224529775023261
$ws.Range("I7").Value = 0.3451699599880819
$ws.Range("J7").Value = 0.3451699599880819
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.03315207120892755
$ws.Range("T7").Value = 0.03315207120892755
$ws.Range("I8").Value = 0.1245448961240849
$ws.Range("J8").Value = 0.1245448961240849
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 83.12110914939734
$ws.Range("R8").Value = 748.089982344576
$ws.Range("S8").Value = 0.03231698544021794
$ws.Range("T8").Value = 0.03231698544021794
$ws.Range("I9").Value = 0.1245448961240849
$ws.Range("J9").Value = 0.1245448961240849
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.08026591588815324
$ws.Range("T9").Value = 0.08026591588815324
$ws.Range("I10").Value = 0.1245448961240849
$ws.Range("J10").Value = 0.1245448961240849
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.01196199479571372
$ws.Range("T10").Value = 0.01196199479571372